# Estadisticos Segundo Parcial 26 Mayo
# On the "Rescatables" sheet, the record for student 23330051920032
# (DE JESUS VERA EDUARDO), previously listed last in this block (row 24),
# is moved to the top of the block (row 16) and its "Reprobadas" count is
# corrected from 1 to 2. The other records (originally rows 16-23) each
# shift down by one row to make room (rows 17-24), keeping their own data
# unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

# Capture the current (pre-edit) values for the block of rows 16-24,
# columns A-G, before overwriting anything.
$firstRow = 16
$lastRow = 24
$numCols = 7

$data = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @()
    for ($c = 1; $c -le $numCols; $c++) {
        $rowVals += , ($ws.Cells.Item($r, $c).Value2)
    }
    $data[$r] = $rowVals
}

# New row 16 = old row 24, with column G (index 7) corrected from 1 to 2.
$movedRow = $data[$lastRow]
$movedRow[6] = 2

$ws.Cells.Item($firstRow, 1).Value = $movedRow[0]
$ws.Cells.Item($firstRow, 2).Value = $movedRow[1]
$ws.Cells.Item($firstRow, 3).Value = $movedRow[2]
$ws.Cells.Item($firstRow, 4).Value = $movedRow[3]
$ws.Cells.Item($firstRow, 5).Value = $movedRow[4]
$ws.Cells.Item($firstRow, 6).Value = $movedRow[5]
$ws.Cells.Item($firstRow, 7).Value = $movedRow[6]

# Old rows 16-23 shift down to become rows 17-24, unchanged otherwise.
for ($r = $firstRow; $r -le ($lastRow - 1); $r++) {
    $destRow = $r + 1
    $srcVals = $data[$r]
    $ws.Cells.Item($destRow, 1).Value = $srcVals[0]
    $ws.Cells.Item($destRow, 2).Value = $srcVals[1]
    $ws.Cells.Item($destRow, 3).Value = $srcVals[2]
    $ws.Cells.Item($destRow, 4).Value = $srcVals[3]
    $ws.Cells.Item($destRow, 5).Value = $srcVals[4]
    $ws.Cells.Item($destRow, 6).Value = $srcVals[5]
    $ws.Cells.Item($destRow, 7).Value = $srcVals[6]
}
